# Auto-generated Excel COM-interop script to apply market-price data updates
# to the Pandaemonium_Profits workbook (columns H-N across ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR sheets), per the scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 101.64
$ws.Range("I33").Value = 57.15
$ws.Range("K33").Value = 57.15
$ws.Range("M33").Value = 171.85
$ws.Range("H63").Value = 38271
$ws.Range("J63").Value = 38271
$ws.Range("L63").Value = 38271
$ws.Range("N63").Value = -39519
$ws.Range("H66").Value = 38271
$ws.Range("J66").Value = 38271
$ws.Range("L66").Value = 114813
$ws.Range("N66").Value = -121053
$ws.Range("H129").Value = 1064.675
$ws.Range("I129").Value = 337
$ws.Range("J129").Value = 1083.3334
$ws.Range("K129").Value = 1011
$ws.Range("L129").Value = 3250.0002
$ws.Range("M129").Value = 3989
$ws.Range("N129").Value = -13250.0002
$ws.Range("H132").Value = 2752.2222
$ws.Range("I132").Value = 2542.6667
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 7628.000100000001
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -5098.000100000001
$ws.Range("N132").Value = -16460

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 70000
$ws.Range("J44").Value = 70000
$ws.Range("L44").Value = 70000
$ws.Range("N44").Value = -70976

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11188.033
$ws.Range("I61").Value = 6484.1113
$ws.Range("J61").Value = 18243.916
$ws.Range("K61").Value = 6484.1113
$ws.Range("L61").Value = 18243.916
$ws.Range("M61").Value = -6272.1113
$ws.Range("N61").Value = -18667.916
$ws.Range("H112").Value = 29257.1
$ws.Range("J112").Value = 29257.1
$ws.Range("L112").Value = 29257.1
$ws.Range("N112").Value = -32211.1
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H120").Value = 79366.336
$ws.Range("J120").Value = 79366.336
$ws.Range("L120").Value = 79366.336
$ws.Range("N120").Value = -89042.336
$ws.Range("H136").Value = 11188.033
$ws.Range("I136").Value = 6484.1113
$ws.Range("J136").Value = 18243.916
$ws.Range("K136").Value = 19452.3339
$ws.Range("L136").Value = 54731.74800000001
$ws.Range("M136").Value = -16902.3339
$ws.Range("N136").Value = -59831.74800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 46120.652
$ws.Range("I134").Value = 2184.0667
$ws.Range("J134").Value = 128501.75
$ws.Range("K134").Value = 6552.2001
$ws.Range("L134").Value = 385505.25
$ws.Range("M134").Value = -4017.2001
$ws.Range("N134").Value = -390575.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2934402.2
$ws.Range("I58").Value = 3953929.8
$ws.Range("J58").Value = 3260.75
$ws.Range("K58").Value = 3953929.8
$ws.Range("L58").Value = 3260.75
$ws.Range("M58").Value = -3953726.8
$ws.Range("N58").Value = -3666.75
$ws.Range("H99").Value = 2195.5
$ws.Range("I99").Value = 1680.75
$ws.Range("J99").Value = 3225
$ws.Range("K99").Value = 1680.75
$ws.Range("L99").Value = 3225
$ws.Range("M99").Value = -182.75
$ws.Range("N99").Value = -6221
$ws.Range("H107").Value = 567.6957
$ws.Range("I107").Value = 311.44446
$ws.Range("J107").Value = 732.4286
$ws.Range("K107").Value = 311.44446
$ws.Range("L107").Value = 732.4286
$ws.Range("M107").Value = 1608.55554
$ws.Range("N107").Value = -4572.4286
$ws.Range("H122").Value = 7600.2
$ws.Range("I122").Value = 8111.3335
$ws.Range("K122").Value = 24334.0005
$ws.Range("M122").Value = -21884.0005
$ws.Range("H126").Value = 2195.5
$ws.Range("I126").Value = 1680.75
$ws.Range("J126").Value = 3225
$ws.Range("K126").Value = 5042.25
$ws.Range("L126").Value = 9675
$ws.Range("M126").Value = -2572.25
$ws.Range("N126").Value = -14615
$ws.Range("H132").Value = 5245.2163
$ws.Range("I132").Value = 5357.9033
$ws.Range("J132").Value = 4663
$ws.Range("K132").Value = 16073.7099
$ws.Range("L132").Value = 13989
$ws.Range("M132").Value = -13543.7099
$ws.Range("N132").Value = -19049
$ws.Range("H134").Value = 3174.2856
$ws.Range("I134").Value = 2611.8096
$ws.Range("K134").Value = 7835.4288
$ws.Range("M134").Value = -5300.4288
$ws.Range("H136").Value = 2934402.2
$ws.Range("I136").Value = 3953929.8
$ws.Range("J136").Value = 3260.75
$ws.Range("K136").Value = 11861789.4
$ws.Range("L136").Value = 9782.25
$ws.Range("M136").Value = -11859239.4
$ws.Range("N136").Value = -14882.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 9804332
$ws.Range("I5").Value = 435.75
$ws.Range("K5").Value = 1307.25
$ws.Range("M5").Value = -1195.25
$ws.Range("H39").Value = 4688.8887
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 4688.8887
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 14066.6661
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -14654.6661
$ws.Range("H80").Value = 2494.6191
$ws.Range("I80").Value = 2568.4285
$ws.Range("J80").Value = 2457.7144
$ws.Range("K80").Value = 7705.2855
$ws.Range("L80").Value = 7373.1432
$ws.Range("M80").Value = -6769.2855
$ws.Range("N80").Value = -9245.143199999999
$ws.Range("H83").Value = 2494.6191
$ws.Range("I83").Value = 2568.4285
$ws.Range("J83").Value = 2457.7144
$ws.Range("K83").Value = 23115.8565
$ws.Range("L83").Value = 22119.4296
$ws.Range("M83").Value = -18435.8565
$ws.Range("N83").Value = -31479.4296
$ws.Range("H122").Value = 864.7619
$ws.Range("I122").Value = 299
$ws.Range("J122").Value = 997.8823
$ws.Range("K122").Value = 2691
$ws.Range("L122").Value = 8980.940699999999
$ws.Range("M122").Value = -241
$ws.Range("N122").Value = -13880.9407
$ws.Range("H133").Value = 6419.1665
$ws.Range("H134").Value = 3984.1035
$ws.Range("I134").Value = 3627.5715
$ws.Range("J134").Value = 4920
$ws.Range("K134").Value = 10882.7145
$ws.Range("L134").Value = 14760
$ws.Range("M134").Value = -5812.7145
$ws.Range("N134").Value = -24900
$ws.Range("H135").Value = 9804332
$ws.Range("I135").Value = 435.75
$ws.Range("K135").Value = 3921.75
$ws.Range("M135").Value = -1386.75
$ws.Range("H137").Value = 31116.264
$ws.Range("I137").Value = 1257.5
$ws.Range("J137").Value = 39078.6
$ws.Range("K137").Value = 3772.5
$ws.Range("L137").Value = 117235.8
$ws.Range("M137").Value = 1327.5
$ws.Range("N137").Value = -127435.8
$ws.Range("H139").Value = 3067.32
$ws.Range("I139").Value = 1790
$ws.Range("J139").Value = 4070.9285
$ws.Range("K139").Value = 5370
$ws.Range("L139").Value = 12212.7855
$ws.Range("M139").Value = -230
$ws.Range("N139").Value = -22492.7855
$ws.Range("H141").Value = 4307.273
$ws.Range("I141").Value = 2913.3333
$ws.Range("K141").Value = 8739.999899999999
$ws.Range("M141").Value = -3559.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 20302
$ws.Range("J103").Value = 20302
$ws.Range("L103").Value = 20302
$ws.Range("N103").Value = -22646
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52216
$ws.Range("H132").Value = 13469.728
$ws.Range("I132").Value = 3777.5
$ws.Range("J132").Value = 15623.556
$ws.Range("K132").Value = 11332.5
$ws.Range("L132").Value = 46870.66800000001
$ws.Range("M132").Value = -8802.5
$ws.Range("N132").Value = -51930.66800000001
$ws.Range("H138").Value = 43459
$ws.Range("J138").Value = 43459
$ws.Range("L138").Value = 43459
$ws.Range("N138").Value = -53739

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 709.9
$ws.Range("I22").Value = 399
$ws.Range("J22").Value = 787.625
$ws.Range("K22").Value = 399
$ws.Range("L22").Value = 787.625
$ws.Range("M22").Value = -104
$ws.Range("N22").Value = -1377.625
$ws.Range("H27").Value = 709.9
$ws.Range("I27").Value = 399
$ws.Range("J27").Value = 787.625
$ws.Range("K27").Value = 399
$ws.Range("L27").Value = 787.625
$ws.Range("M27").Value = -292
$ws.Range("N27").Value = -1001.625
$ws.Range("H62").Value = 38249
$ws.Range("J62").Value = 38249
$ws.Range("L62").Value = 38249
$ws.Range("N62").Value = -39497
$ws.Range("H64").Value = 36150
$ws.Range("J64").Value = 36150
$ws.Range("L64").Value = 36150
$ws.Range("N64").Value = -36600
$ws.Range("H65").Value = 38249
$ws.Range("J65").Value = 38249
$ws.Range("L65").Value = 114747
$ws.Range("N65").Value = -120987
$ws.Range("H67").Value = 36150
$ws.Range("J67").Value = 36150
$ws.Range("L67").Value = 36150
$ws.Range("N67").Value = -37710

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 2401
$ws.Range("I8").Value = 2401
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 2401
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -2261
$ws.Range("N8").ClearContents()
$ws.Range("H39").Value = 30044
$ws.Range("I39").Value = 30044
$ws.Range("K39").Value = 30044
$ws.Range("M39").Value = -29631
$ws.Range("H100").Value = 1176.3846
$ws.Range("I100").Value = 456.42856
$ws.Range("J100").Value = 2016.3334
$ws.Range("K100").Value = 912.85712
$ws.Range("L100").Value = 4032.6668
$ws.Range("M100").Value = -371.85712
$ws.Range("N100").Value = -5114.6668
$ws.Range("H124").Value = 25031.6
$ws.Range("J124").Value = 25031.6
$ws.Range("L124").Value = 25031.6
$ws.Range("N124").Value = -34851.6
